$d = $word.ActiveDocument

# PiOTR paragraph: merge the "collectables behaviours" run with the following
# ",   car spin " run into a single run (search only within the second run's
# text so the preceding space-run is left untouched, matching the target XML).
$d.Content.Find.Execute(",   car spin ", $true, $false, $false, $false, $false, `
    $true, 1, $false, ",   car spin ", 2)

# James paragraph: merge the "collectables behaviours" run with the following
# ", fix fuel. Change " run into a single run, same technique.
$d.Content.Find.Execute(", fix fuel. Change ", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", fix fuel. Change ", 2)

# Mike's line: drop "and fix background animation."
$d.Content.Find.Execute("Mike- add level design and fix background animation. ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Mike- add level design ", 2)
